$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "prem1227"
$ws.Range("C2").Value = "shankar1227"
$ws.Range("B2").Value = "shankar77@gmail.com"
$ws.Range("B3").Value = "prem28@gmail.com"

$ws.Range("B3").Select()
